$d = $word.ActiveDocument

# 1) Replace the placeholder text in the first paragraph, and drop the trailing space run.
$d.Content.Find.Execute("**ID__AFFARS_mp_5332_7_topic_1__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP5332_7__ID**", 2)

# 2) Update the first paragraph's indentation and add a paragraph border.
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 11.25
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
